# Fruta / hortaliza, semanal
# Insert two new daily-price rows (for 2021-09-20 = serial 44463) into the
# "Naranja" sheet right after the current header + first data block, pushing
# the remaining rows down, and append the freed-up rows at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 224:282 down to 226:284, freeing up rows 224 and 225 for the
# two new records.
$ws.Rows("224:225").Insert()

# New row 224: Naranja / Lane Late / Primera
$ws.Range("A224").Value2 = 5
$ws.Range("B224").Value2 = "Macroferia Regional de Talca"
$ws.Range("C224").Value2 = "Maule"
$ws.Range("D224").Value2 = 44463
$ws.Range("E224").Value2 = 7
$ws.Range("F224").Value2 = "Fruta"
$ws.Range("G224").Value2 = 100102
$ws.Range("H224").Value2 = "Cítricos"
$ws.Range("I224").Value2 = 100102005
$ws.Range("J224").Value2 = "Naranja"
$ws.Range("K224").Value2 = "Lane Late"
$ws.Range("L224").Value2 = "Primera"
$ws.Range("M224").Value2 = 350
$ws.Range("N224").Value2 = 6000
$ws.Range("O224").Value2 = 6000
$ws.Range("P224").Value2 = 6000
$ws.Range("Q224").Value2 = "`$/bandeja 15 kilos granel"
$ws.Range("R224").Value2 = "Región de O'Higgins"
$ws.Range("S224").Value2 = 400
$ws.Range("T224").Value2 = 15

# New row 225: Naranja / Navel Late / Primera
$ws.Range("A225").Value2 = 5
$ws.Range("B225").Value2 = "Macroferia Regional de Talca"
$ws.Range("C225").Value2 = "Maule"
$ws.Range("D225").Value2 = 44463
$ws.Range("E225").Value2 = 7
$ws.Range("F225").Value2 = "Fruta"
$ws.Range("G225").Value2 = 100102
$ws.Range("H225").Value2 = "Cítricos"
$ws.Range("I225").Value2 = 100102005
$ws.Range("J225").Value2 = "Naranja"
$ws.Range("K225").Value2 = "Navel Late"
$ws.Range("L225").Value2 = "Primera"
$ws.Range("M225").Value2 = 300
$ws.Range("N225").Value2 = 7000
$ws.Range("O225").Value2 = 7000
$ws.Range("P225").Value2 = 7000
$ws.Range("Q225").Value2 = "`$/bandeja 15 kilos granel"
$ws.Range("R225").Value2 = "Región de O'Higgins"
$ws.Range("S225").Value2 = 467
$ws.Range("T225").Value2 = 15
